# "Generate Report for Handoff"
#
# This updates the localization-status report:
#  - Marks the b.md row as "Ready for handoff" (status) with a fresh
#    "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamp
#    on every sheet.
#  - Records the newly generated handoff xliff file names for b.md on the
#    zh-cn and de-de sheets.
#  - Flags the zh-cn b.md row's "Content Duplicate" as False (it is a new
#    file) and records an "Error Detail" message (on zh-cn and de-de) that
#    the handback file is stale.
#  - Widens the "Error Detail" column now that it holds a long message.

$wb = $excel.ActiveWorkbook

# Helper: assign literal text to a cell without Excel's auto-coercion of
# recognisable literals ("True"/"False" -> Boolean, etc.). We stash the
# text as a formula result in a scratch cell, copy it, and paste-special
# just the value back onto the destination - this keeps the destination
# a plain shared-string cell (no formula, no style change) the same way
# a direct text write would.
function Set-LiteralText($range, [string]$text) {
    $sheet = $range.Worksheet
    $scratch = $sheet.Range("ZZ1")
    $escaped = $text.Replace("""", """""")
    $scratch.Formula = "=""" + $escaped + """"
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md file.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-25 18:39:18"

# ---------------------------------------------------------------------
# zh-cn sheet: row 2 is a.md, row 3 is b.md.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Ready for handoff"

$zh.Range("C3").Value = "Ready for handoff"
Set-LiteralText $zh.Range("F3") "False"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-25 18:39:14"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71644b2cdc10670041e8d5993ba05b4b8506e79c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d4e749be65f7dfd5d4f59f4d2af614ae92d7e7b/e2e/b.md."

$zh.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet: row 2 is a.md, row 3 is b.md.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Ready for handoff"

$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-08-25 18:39:18"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71644b2cdc10670041e8d5993ba05b4b8506e79c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d4e749be65f7dfd5d4f59f4d2af614ae92d7e7b/e2e/b.md."

$de.Columns.Item(16).ColumnWidth = 39.14
